$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.262.55"
$ws.Range("E2").Value = "  -5.33%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.669.47"
$ws.Range("E3").Value = "  -3.12%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.006"
$ws.Range("E4").Value = "  +0.38%  "

# Row 5 - BNB
Set-TextValue "D5" "217.73"
$ws.Range("E5").Value = "  -3.52%  "

# Row 6 - XRP
Set-TextValue "D6" "0.5087"
$ws.Range("E6").Value = "  -11.27%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.28%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.2659"
$ws.Range("E8").Value = "  -2.03%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.06373"
$ws.Range("E9").Value = "  -3.22%  "

# Row 10 - Solana
Set-TextValue "D10" "21.45"
$ws.Range("E10").Value = "  -6.06%  "

# Row 11 - TRON
Set-TextValue "D11" "0.07373"
$ws.Range("E11").Value = "  -2.08%  "

# Row 12 - WrappedEther
Set-TextValue "D12" "1.676.73"
$ws.Range("E12").Value = "  -2.90%  "

# Row 13 - Polkadot
Set-TextValue "D13" "4.548"
$ws.Range("E13").Value = "  -2.60%  "

# Row 14 - Polygon
Set-TextValue "D14" "0.5816"
$ws.Range("E14").Value = "  -2.70%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "1.897.83"
$ws.Range("E15").Value = "  -3.06%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.000008531"
$ws.Range("E16").Value = "  -0.88%  "

# Row 17 - Litecoin
Set-TextValue "D17" "64.67"
$ws.Range("E17").Value = "  -12.82%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "26.342.94"
$ws.Range("E18").Value = "  -4.84%  "

# Row 19 - Uniswap
Set-TextValue "D19" "4.940"

# Row 20 - Dai
Set-TextValue "D20" "1.006"
$ws.Range("E20").Value = "  +0.15%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "189.64"
$ws.Range("E22").Value = "  -6.95%  "

# Row 23 - Chainlink
Set-TextValue "D23" "6.203"
$ws.Range("E23").Value = "  -5.48%  "

# Row 24 - BinanceUSD
Set-TextValue "D24" "1.008"
$ws.Range("E24").Value = "  +0.36%  "

# Row 25 - Monero
Set-TextValue "D25" "143.60"
$ws.Range("E25").Value = "  -4.03%  "

# Row 26 - Cosmos
Set-TextValue "D26" "7.648"
$ws.Range("E26").Value = "  -4.48%  "

# Row 27 - Stellar
Set-TextValue "D27" "0.1176"
$ws.Range("E27").Value = "  -3.64%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "15.63"
$ws.Range("E28").Value = "  -2.93%  "

# Row 29 - Hedera
Set-TextValue "D29" "0.05865"
$ws.Range("E29").Value = "  -4.67%  "

# Row 30 - Toncoin
Set-TextValue "D30" "1.272"
$ws.Range("E30").Value = "  -7.76%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.321"
$ws.Range("E31").Value = "  -4.86%  "

# Row 32 - now Filecoin (was InternetComputer(DFINITY))
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "3.513"
$ws.Range("E32").Value = "  -5.21%  "

# Row 33 - now InternetComputer(DFINITY) (was Filecoin)
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D33" "3.522"
$ws.Range("E33").Value = "  -4.86%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -1.88%  "

# Row 35 - ARBITRUM
Set-TextValue "D35" "1.012"
$ws.Range("E35").Value = "  -1.63%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "0.5998"
$ws.Range("E36").Value = "  -6.50%  "

# Row 37 - HuobiToken
Set-TextValue "D37" "2.361"
$ws.Range("E37").Value = "  -2.79%  "

# Row 38 - MXToken
Set-TextValue "D38" "2.646"
$ws.Range("E38").Value = "  -1.17%  "

# Row 39 - VeChain
Set-TextValue "D39" "0.01615"
$ws.Range("E39").Value = "  -2.60%  "

# Row 40 - FraxShare
Set-TextValue "D40" "6.035"

# Row 41 - Maker
Set-TextValue "D41" "1.074.54"
$ws.Range("E41").Value = "  -3.90%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "0.8680"
$ws.Range("E42").Value = "  -0.77%  "

# Row 43 - PaxDollar
$ws.Range("E43").Value = "  +0.50%  "

# Row 44 - Quant
Set-TextValue "D44" "99.70"
$ws.Range("E44").Value = "  +0.28%  "

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.817.17"
$ws.Range("E45").Value = "  -2.96%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  +1.38%  "

# Row 47 - Aave
Set-TextValue "D47" "55.85"
$ws.Range("E47").Value = "  -5.22%  "

# Row 48 - Frax
$ws.Range("E48").Value = "  +0.60%  "

# Row 49 - EnergySwap
Set-TextValue "D49" "8.062"
$ws.Range("E49").Value = "  -1.55%  "

# Row 50 - Mantle
Set-TextValue "D50" "0.4299"
$ws.Range("E50").Value = "  -2.50%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -3.40%  "
